# Update the EJ43 sheet: new column headers (J4301-J4306) and fill in the
# previously empty measurement tables, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ43")

# --- Header row (B1:G1) ---------------------------------------------------
$ws.Range("B1").Value = "J4301"
$ws.Range("C1").Value = "J4302"
$ws.Range("D1").Value = "J4303"
$ws.Range("E1").Value = "J4304"
$ws.Range("F1").Value = "J4305"
$ws.Range("G1").Value = "J4306"

# --- Data rows 2-8 (dx = 0, 0.5, 1, 1.5, 2, 2.5, 3) ------------------------
$ws.Range("B2").Value = 298.60000000000002
$ws.Range("C2").Value = 221.8
$ws.Range("D2").Value = 173.8
$ws.Range("E2").Value = 133.69999999999999
$ws.Range("F2").Value = 98.4
$ws.Range("G2").Value = 71.5

$ws.Range("B3").Value = 309.2
$ws.Range("C3").Value = 229.3
$ws.Range("D3").Value = 183.5
$ws.Range("E3").Value = 136.5
$ws.Range("F3").Value = 100.8
$ws.Range("G3").Value = 73.5

$ws.Range("B4").Value = 320.39999999999998
$ws.Range("C4").Value = 238.7
$ws.Range("D4").Value = 192.2
$ws.Range("E4").Value = 139.19999999999999
$ws.Range("F4").Value = 102.9
$ws.Range("G4").Value = 74.8

$ws.Range("B5").Value = 330.3
$ws.Range("C5").Value = 248.2
$ws.Range("D5").Value = 201.1
$ws.Range("E5").Value = 144.80000000000001
$ws.Range("F5").Value = 105
$ws.Range("G5").Value = 76.7

$ws.Range("B6").Value = 340.5
$ws.Range("C6").Value = 257.2
$ws.Range("D6").Value = 209.7
$ws.Range("E6").Value = 148.69999999999999
$ws.Range("F6").Value = 107.5
$ws.Range("G6").Value = 79.400000000000006

$ws.Range("B7").Value = 350.5
$ws.Range("C7").Value = 265.2
$ws.Range("D7").Value = 218.2
$ws.Range("E7").Value = 151.9
$ws.Range("F7").Value = 109.8
$ws.Range("G7").Value = 81.3

$ws.Range("B8").Value = 359.2
$ws.Range("C8").Value = 271.7
$ws.Range("D8").Value = 226.9
$ws.Range("E8").Value = 155.30000000000001
$ws.Range("F8").Value = 112.7
$ws.Range("G8").Value = 82.6

# --- Move the active selection from B2 to C2 -------------------------------
$ws.Activate()
$ws.Range("C2").Select()
